# Reverse the order of the comma-separated entries in the "Recorded By"
# column (column G) for every data row in the active worksheet.
#
# e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"
#      "System, backup@backdoor.com, system" -> "system, backup@backdoor.com, System"
# Cells that contain only a single value (no comma) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

# Find the "Recorded By" column dynamically by scanning the header row (row 1)
$headerRow = 1
$recordedByCol = 0
$lastCol = $usedRange.Columns.Count + $usedRange.Column - 1
for ($c = 1; $c -le $lastCol; $c++) {
    $headerValue = $ws.Cells.Item($headerRow, $c).Value2
    if ($headerValue -eq "Recorded By") {
        $recordedByCol = $c
        break
    }
}

if ($recordedByCol -eq 0) {
    $recordedByCol = 7   # fall back to column G
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $recordedByCol)
    $value = $cell.Value2

    if ($value -ne $null -and $value -ne "") {
        $parts = $value -split ","
        if ($parts.Count -gt 1) {
            $trimmed = @()
            foreach ($p in $parts) {
                $trimmed += $p.Trim()
            }

            $reversed = @()
            for ($i = $trimmed.Count - 1; $i -ge 0; $i--) {
                $reversed += $trimmed[$i]
            }

            $newValue = [string]::Join(", ", $reversed)
            $cell.Value2 = $newValue
        }
    }
}
